# SwaadSutra_Daily_2026-01-26.xlsx - add new order (#30, Minakshi) placed above
# the existing order (#29, Priyanka Patil), and update the Summary + Items
# Breakdown sheets to reflect the new order.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Daily Orders" -------------------------------------------------
# Insert a brand-new row above the current row 2 so the new order (#30)
# becomes the first data row, and the previous order (#29) shifts to row 3.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = 30
$ws1.Range("B2").Value = "'2026-01-26 12:02"
$ws1.Range("C2").Value = "Minakshi"
$ws1.Range("D2").Value = "A201"
$ws1.Range("E2").Value = "'7387851735"
$ws1.Range("F2").Value = "Wheat Chapati x5"
$ws1.Range("G2").Value = 75
$ws1.Range("H2").Value = "NEW"
$ws1.Range("I2").Value = "PENDING"
$ws1.Range("J2").Value = "'2026-01-26"
$ws1.Range("K2").Value = "'19:30"
$ws1.Range("L2").Value = "'"
$ws1.Range("M2").Value = "'"
$ws1.Range("N2").Value = "'"

# --- Sheet 2: "Summary" ------------------------------------------------------
# One more order was placed, and it is also a "New" order; the total revenue
# grows by the new order's total (75).
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = 2
$ws2.Range("G2").Value = 115

# --- Sheet 3: "Items Breakdown" ----------------------------------------------
# Insert a new row above the existing "Jawar Bhakari" row for the newly
# ordered item "Wheat Chapati" (qty 5, revenue 75).
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()

$ws3.Range("A2").Value = "Wheat Chapati"
$ws3.Range("B2").Value = 5
$ws3.Range("C2").Value = 75
